# Auto-generated Excel COM-interop script
# Applies numeric corrections to the pricing columns (H-N) across several
# sheets in the workbook, matching the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 271.33334
$ws.Range("I8").Value = 271.33334
$ws.Range("K8").Value = 814.0000200000001
$ws.Range("M8").Value = -675.0000200000001

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H17").Value = 2938.5
$ws.Range("J17").Value = 2938.5
$ws.Range("L17").Value = 8815.5
$ws.Range("N17").Value = -9151.5

$ws.Range("H28").Value = 3999.5
$ws.Range("I28").Value = 3333
$ws.Range("K28").Value = 3333
$ws.Range("M28").Value = -2848

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H51").Value = 10316
$ws.Range("I51").Value = 9047.333000000001
$ws.Range("K51").Value = 9047.333000000001
$ws.Range("M51").Value = -8563.333000000001

$ws.Range("H55").Value = 1842.8889
$ws.Range("I55").Value = 212
$ws.Range("J55").Value = 2658.3333
$ws.Range("K55").Value = 212
$ws.Range("L55").Value = 2658.3333
$ws.Range("M55").Value = 2
$ws.Range("N55").Value = -3086.3333

$ws.Range("H70").Value = 4655.5557
$ws.Range("I70").Value = 4050
$ws.Range("J70").Value = 9500
$ws.Range("K70").Value = 12150
$ws.Range("L70").Value = 28500
$ws.Range("M70").Value = -11880
$ws.Range("N70").Value = -29040

$ws.Range("H73").Value = 4655.5557
$ws.Range("I73").Value = 4050
$ws.Range("J73").Value = 9500
$ws.Range("K73").Value = 12150
$ws.Range("L73").Value = 28500
$ws.Range("M73").Value = -11214
$ws.Range("N73").Value = -30372

$ws.Range("H100").Value = 2475.9092
$ws.Range("I100").Value = 2555.8096
$ws.Range("J100").Value = 798
$ws.Range("K100").Value = 2555.8096
$ws.Range("L100").Value = 798
$ws.Range("M100").Value = -2014.8096
$ws.Range("N100").Value = -1880

$ws.Range("H107").Value = 599.5
$ws.Range("I107").Value = 599
$ws.Range("K107").Value = 599
$ws.Range("M107").Value = 1321

$ws.Range("H137").Value = 3124.5557
$ws.Range("I137").Value = 2928.3333
$ws.Range("J137").Value = 3222.6667
$ws.Range("K137").Value = 8784.999899999999
$ws.Range("L137").Value = 9668.000100000001
$ws.Range("M137").Value = -6234.999899999999
$ws.Range("N137").Value = -14768.0001

$ws.Range("H141").Value = 1850
$ws.Range("I141").Value = 1850
$ws.Range("K141").Value = 5550
$ws.Range("M141").Value = -370

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3368.56
$ws.Range("I32").Value = 2963.625
$ws.Range("K32").Value = 2963.625
$ws.Range("M32").Value = -2676.625

$ws.Range("H36").Value = 2066.25
$ws.Range("I36").Value = 2066.25
$ws.Range("K36").Value = 2066.25
$ws.Range("M36").Value = -1720.25

$ws.Range("H45").Value = 1628.3158
$ws.Range("I45").Value = 1307.6666
$ws.Range("K45").Value = 1307.6666
$ws.Range("M45").Value = -930.6666

$ws.Range("H74").Value = 4312.7856
$ws.Range("J74").Value = 4985.8
$ws.Range("L74").Value = 4985.8
$ws.Range("N74").Value = -6733.8

$ws.Range("H77").Value = 4312.7856
$ws.Range("J77").Value = 4985.8
$ws.Range("L77").Value = 24929
$ws.Range("N77").Value = -33665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 4649.8
$ws.Range("J64").Value = 4649.8
$ws.Range("L64").Value = 4649.8
$ws.Range("N64").Value = -5099.8

$ws.Range("H67").Value = 4649.8
$ws.Range("J67").Value = 4649.8
$ws.Range("L67").Value = 4649.8
$ws.Range("N67").Value = -6209.8

$ws.Range("H86").Value = 9705.823
$ws.Range("I86").Value = 15168.25
$ws.Range("J86").Value = 4850.3335
$ws.Range("K86").Value = 15168.25
$ws.Range("L86").Value = 4850.3335
$ws.Range("M86").Value = -14045.25
$ws.Range("N86").Value = -7096.3335

$ws.Range("H89").Value = 9705.823
$ws.Range("I89").Value = 15168.25
$ws.Range("J89").Value = 4850.3335
$ws.Range("K89").Value = 75841.25
$ws.Range("L89").Value = 24251.6675
$ws.Range("M89").Value = -70225.25
$ws.Range("N89").Value = -35483.6675

$ws.Range("H94").Value = 4814.1333
$ws.Range("I94").Value = 4907.1665
$ws.Range("J94").Value = 4752.1113
$ws.Range("K94").Value = 4907.1665
$ws.Range("L94").Value = 4752.1113
$ws.Range("M94").Value = -4456.1665
$ws.Range("N94").Value = -5654.1113

$ws.Range("H96").Value = 24994.5
$ws.Range("I96").Value = 24994.5
$ws.Range("K96").Value = 24994.5
$ws.Range("M96").Value = -22248.5

$ws.Range("H105").Value = 3274.5
$ws.Range("I105").Value = 2739.4
$ws.Range("K105").Value = 2739.4
$ws.Range("M105").Value = -992.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10574

$ws.Range("H31").Value = 18970.75
$ws.Range("I31").Value = 24296.334
$ws.Range("K31").Value = 24296.334
$ws.Range("M31").Value = -24001.334

$ws.Range("H32").Value = 14500
$ws.Range("J32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("N32").Value = -25632

$ws.Range("H34").Value = 18970.75
$ws.Range("I34").Value = 24296.334
$ws.Range("K34").Value = 24296.334
$ws.Range("M34").Value = -24094.334

$ws.Range("H99").Value = 1453.3636
$ws.Range("I99").Value = 1327.7142
$ws.Range("J99").Value = 1673.25
$ws.Range("K99").Value = 1327.7142
$ws.Range("L99").Value = 1673.25
$ws.Range("M99").Value = 170.2858000000001
$ws.Range("N99").Value = -4669.25

$ws.Range("H105").Value = 1038.2727
$ws.Range("I105").Value = 699.55554
$ws.Range("K105").Value = 699.55554
$ws.Range("M105").Value = 1047.44446

$ws.Range("H126").Value = 1453.3636
$ws.Range("I126").Value = 1327.7142
$ws.Range("J126").Value = 1673.25
$ws.Range("K126").Value = 3983.1426
$ws.Range("L126").Value = 5019.75
$ws.Range("M126").Value = -1513.1426
$ws.Range("N126").Value = -9959.75

$ws.Range("H132").Value = 4804.3335
$ws.Range("I132").Value = 3877.4546
$ws.Range("K132").Value = 11632.3638
$ws.Range("M132").Value = -9102.363799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 359
$ws.Range("J41").Value = 359
$ws.Range("L41").Value = 1077
$ws.Range("N41").Value = -1753

$ws.Range("H43").Value = 12445
$ws.Range("J43").Value = 12445
$ws.Range("L43").Value = 37335
$ws.Range("N43").Value = -37563

$ws.Range("H68").Value = 273.66666
$ws.Range("I68").Value = 273.66666
$ws.Range("K68").Value = 820.9999799999999
$ws.Range("M68").Value = -9.999979999999937

$ws.Range("H71").Value = 273.66666
$ws.Range("I71").Value = 273.66666
$ws.Range("K71").Value = 2462.99994
$ws.Range("M71").Value = 1593.00006

$ws.Range("H113").Value = 1040.6
$ws.Range("I113").Value = 1045.174
$ws.Range("K113").Value = 3135.522
$ws.Range("M113").Value = -965.5219999999999

$ws.Range("H131").Value = 2779.6667
$ws.Range("J131").Value = 3549.6667
$ws.Range("L131").Value = 10649.0001
$ws.Range("N131").Value = -20729.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 80
$ws.Range("I2").Value = 80
$ws.Range("K2").Value = 80
$ws.Range("M2").Value = 32

$ws.Range("H22").Value = 866.8570999999999
$ws.Range("I22").Value = 674.5
$ws.Range("J22").Value = 1123.3334
$ws.Range("K22").Value = 674.5
$ws.Range("L22").Value = 1123.3334
$ws.Range("M22").Value = -379.5
$ws.Range("N22").Value = -1713.3334

$ws.Range("H27").Value = 866.8570999999999
$ws.Range("I27").Value = 674.5
$ws.Range("J27").Value = 1123.3334
$ws.Range("K27").Value = 674.5
$ws.Range("L27").Value = 1123.3334
$ws.Range("M27").Value = -567.5
$ws.Range("N27").Value = -1337.3334

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H100").Value = 1594.3334
$ws.Range("I100").Value = 1594.3334
$ws.Range("K100").Value = 1594.3334
$ws.Range("M100").Value = -1053.3334

$ws.Range("H132").Value = 4447.875
$ws.Range("I132").Value = 4098.1665
$ws.Range("K132").Value = 12294.4995
$ws.Range("M132").Value = -9764.499500000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5832.6665
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 7249.5
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 7249.5
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -8497.5

$ws.Range("H65").Value = 5832.6665
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 7249.5
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 36247.5
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -42487.5

$ws.Range("H100").Value = 281.8
$ws.Range("I100").Value = 227.25
$ws.Range("K100").Value = 454.5
$ws.Range("M100").Value = 86.5

$ws.Range("H126").Value = 1402
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 1695.8334
$ws.Range("I132").Value = 1695.8334
$ws.Range("K132").Value = 5087.5002
$ws.Range("M132").Value = -2557.5002

$ws.Range("H136").Value = 12250
$ws.Range("J136").Value = 21296.334
$ws.Range("L136").Value = 63889.00199999999
$ws.Range("N136").Value = -68989.00199999999
